# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.458.84"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "3.465.25"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.78"
$ws.Range("E5").Value = "  +0.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.56"
$ws.Range("E6").Value = "  -2.29%  "

$ws.Range("D7").Value = "3.462.22"
$ws.Range("E7").Value = "  +0.55%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  -1.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  -1.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.87"
$ws.Range("E11").Value = "  +6.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.413"
$ws.Range("E12").Value = "  -1.98%  "

$ws.Range("D13").Value = "4.054.38"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000210"
$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.64"
$ws.Range("E15").Value = "  -2.32%  "

$ws.Range("D16").Value = "67.190.51"
$ws.Range("E16").Value = "  +0.70%  "

$ws.Range("D17").Value = "3.465.08"
$ws.Range("E17").Value = "  -0.74%  "

$ws.Range("E18").Value = "  +0.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.60"
$ws.Range("E19").Value = "  +7.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.22"
$ws.Range("E20").Value = "  -2.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.18"
$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "424.74"
$ws.Range("E22").Value = "  -2.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.596"
$ws.Range("E23").Value = "  -2.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.73"
$ws.Range("E24").Value = "  +0.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"

$ws.Range("D26").Value = "3.603.03"
$ws.Range("E26").Value = "  +0.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000114"
$ws.Range("E27").Value = "  -3.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.59"
$ws.Range("E28").Value = "  -1.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.02"
$ws.Range("E29").Value = "  -3.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.49"
$ws.Range("E30").Value = "  +1.25%  "

$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.52"
$ws.Range("E32").Value = "  -4.40%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.164"
$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.12"
$ws.Range("E34").Value = "  -0.49%  "

$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.74"
$ws.Range("E36").Value = "  -2.60%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.65"
$ws.Range("E37").Value = "  -6.86%  "

$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.84"
$ws.Range("E38").Value = "  -0.06%  "

$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "174.23"
$ws.Range("E40").Value = "  +0.40%  "

$ws.Range("E41").Value = "  +0.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.23"
$ws.Range("E42").Value = "  -1.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.887"
$ws.Range("E43").Value = "  +0.89%  "

$ws.Range("E44").Value = "  -10.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "45.91"
$ws.Range("E45").Value = "  -0.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "27.36"
$ws.Range("E46").Value = "  -6.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.19"
$ws.Range("E47").Value = "  -3.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.25"
$ws.Range("E48").Value = "  -2.84%  "

$ws.Range("B49").Value = "SuiNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.961"

$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.33"
$ws.Range("E50").Value = "  -3.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.241"
$ws.Range("E51").Value = "  -0.75%  "

